$d = $word.ActiveDocument
$bullet = [char]8226

# ---------------------------------------------------------------------------
# Change 1: Collapse the three detailed "CORE COMPETENCIES" paragraphs into a
# single short summary paragraph.
# ---------------------------------------------------------------------------

$p6 = $d.Paragraphs.Item(6)
$p7 = $d.Paragraphs.Item(7)
$p8 = $d.Paragraphs.Item(8)

if ($p6.Range.Text -notmatch "^Statistical Analysis") {
    throw "Unexpected paragraph 6 text: $($p6.Range.Text)"
}
if ($p7.Range.Text -notmatch "^Big Data") {
    throw "Unexpected paragraph 7 text: $($p7.Range.Text)"
}
if ($p8.Range.Text -notmatch "^Data Visualization") {
    throw "Unexpected paragraph 8 text: $($p8.Range.Text)"
}

# Delete paragraphs 7 and 8 entirely (their marks merge away); paragraph 6's
# own mark survives so it keeps its original paragraph formatting.
$rngRemove = $d.Range($p7.Range.Start, $p8.Range.End)
$rngRemove.Delete()

# Replace paragraph 6's run text with the condensed summary line.
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Text = "Statistical Analysis & Machine Learning $bullet Big Data & Data Engineering $bullet Data Visualization & Reporting"

# ---------------------------------------------------------------------------
# Change 2: Insert a new "TECHNICAL SKILLS" section (one Heading2 paragraph
# plus three body paragraphs) right after the "Led multi-million dollar ..."
# bullet and before the closing "For a more detailed ..." paragraph.
# ---------------------------------------------------------------------------

$ledIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "Led multi-million dollar research projects") {
        $ledIndex = $i
        break
    }
}
if ($ledIndex -eq 0) {
    throw "Could not find the 'Led multi-million dollar...' paragraph"
}

$newTexts = @(
    "TECHNICAL SKILLS",
    "STATISTICAL ANALYSIS & MACHINE LEARNING Advanced Statistical Modeling; Predictive Analytics; Data Mining; Machine Learning",
    "BIG DATA & DATA ENGINEERING Big Data Processing; Data Warehousing; Cloud Platforms; Data Pipeline Optimization",
    "DATA VISUALIZATION & REPORTING Data Visualization; Geospatial Analysis; Interactive Dashboards; Business Intelligence"
)

$cur = $ledIndex
foreach ($text in $newTexts) {
    $p = $d.Paragraphs.Item($cur)
    $null = $p.Range.InsertParagraphAfter()
    $cur = $cur + 1
    $newP = $d.Paragraphs.Item($cur)
    $newP.Range.Text = $text
}

# The first of the newly-inserted paragraphs ("TECHNICAL SKILLS") becomes a
# Heading 2, matching the other section headings in the document.
$headingPara = $d.Paragraphs.Item($ledIndex + 1)
$headingPara.Style = "Heading 2"
